# Automatische test-sync: 2025-06-19 21:21:50
# Appends the new "Vragen over samenwerking" mail-log row to the Logs sheet
# and refreshes the Dashboard category-count table so the
# "Samenwerking / Partnerverzoek" row reflects the new total and is
# re-sorted ahead of "Productinformatie".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs: new row 12 -------------------------------------------------
$logs.Range("A12").Value = "Vragen over samenwerking"
$logs.Range("B12").Value = "mailmind.test@zohomail.eu"
$logs.Range("C12").Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Range("D12").Value = "Samenwerking / Partnerverzoek"
$logs.Range("F12").Value = "2025-06-19 21:21:17"
$logs.Range("G12").Value = "Nee"

# --- Logs: grow the conditional-formatting ranges to include row 12 ----
$catRules = $logs.Range("D2:D11").FormatConditions
for ($i = 1; $i -le $catRules.Count; $i++) {
    $catRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D12"))
}

$answeredRules = $logs.Range("G2:G11").FormatConditions
for ($i = 1; $i -le $answeredRules.Count; $i++) {
    $answeredRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G12"))
}

# --- Dashboard: update counts & re-sort --------------------------------
# "Samenwerking / Partnerverzoek" now has 3 hits (was 2); it moves above
# "Productinformatie" (still 3) since it now ties and takes the new mail's
# category precedence.
$dashboard.Range("A2").Value = "Samenwerking / Partnerverzoek"
$dashboard.Range("B2").Value = 3
$dashboard.Range("A3").Value = "Productinformatie"
$dashboard.Range("B3").Value = 3
